$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C: plain numeric "Qty executed upto date" values
$ws.Range("C8").Value = 29
$ws.Range("C9").Value = 95
$ws.Range("C10").Value = 85
$ws.Range("C11").Value = 63
$ws.Range("C12").Value = 34
$ws.Range("C13").Value = 90
$ws.Range("C14").Value = 18
$ws.Range("C15").Value = 84
$ws.Range("C16").Value = 34
$ws.Range("C17").Value = 40

# Column G/H cells hold amounts stored as text (2-decimal formatted strings).
# Force text entry (so the "24320.00"-style value isn't reparsed as a plain
# number), then clear the temporary formatting so no stray style survives.
$textCells = "G9","G10","G11","G13","G14","G19","H19","G21","H21"
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("G9").Value = "24320.00"
$ws.Range("G10").Value = "40120.00"
$ws.Range("G11").Value = "41706.00"
$ws.Range("G13").Value = "12240.00"
$ws.Range("G14").Value = "414.00"
$ws.Range("G19").Value = "118800.00"
$ws.Range("H19").Value = "118800.00"
$ws.Range("G21").Value = "118800.00"
$ws.Range("H21").Value = "118800.00"

foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats()
}
